$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.64
$ws.Range("G2").Value = 3.75
$ws.Range("H2").Value = 2.26
$ws.Range("I2").Value = 2.98
$ws.Range("J2").Value = 3.2
$ws.Range("K2").Value = 5.6
$ws.Range("F3").Value = 1.86
$ws.Range("G3").Value = 1.87
$ws.Range("I3").Value = 4.7
$ws.Range("K3").Value = 4.1
$ws.Range("N3").Value = 5.4
$ws.Range("O3").Value = 1.2
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 1.63
$ws.Range("R3").Value = 1.59
$ws.Range("S3").Value = 2.56
$ws.Range("T3").Value = 1.61
$ws.Range("U3").Value = 2.52
$ws.Range("X3").Value = 22
$ws.Range("Y3").Value = 23
$ws.Range("Z3").Value = 40
$ws.Range("AA3").Value = 110
$ws.Range("AB3").Value = 13
$ws.Range("AC3").Value = 9.6
$ws.Range("AD3").Value = 19
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 14.5
$ws.Range("AI3").Value = 1000
$ws.Range("AM3").Value = 65
$ws.Range("AN3").Value = 8.800000000000001
$ws.Range("AO3").Value = 40
$ws.Range("F4").Value = 7.4
$ws.Range("G4").Value = 14.5
$ws.Range("H4").Value = 1.32
$ws.Range("I4").Value = 1.39
$ws.Range("K4").Value = 6.4
$ws.Range("N4").Value = 4.8
$ws.Range("O4").Value = 1.2
$ws.Range("P4").Value = 2.32
$ws.Range("Q4").Value = 1.54
$ws.Range("S4").Value = 2.32
$ws.Range("T4").Value = 1.98
$ws.Range("U4").Value = 1.8
$ws.Range("X4").Value = 29
$ws.Range("Y4").Value = 11
$ws.Range("Z4").Value = 10.5
$ws.Range("AA4").Value = 13
$ws.Range("AB4").Value = 44
$ws.Range("AC4").Value = 16
$ws.Range("AD4").Value = 13
$ws.Range("AE4").Value = 18
$ws.Range("AF4").Value = 130
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 36
$ws.Range("AI4").Value = 46
$ws.Range("AJ4").Value = 500
$ws.Range("AK4").Value = 220
$ws.Range("AL4").Value = 170
$ws.Range("AM4").Value = 190
$ws.Range("AN4").Value = 280
$ws.Range("AO4").Value = 6.2
$ws.Range("H5").Value = 3.95
$ws.Range("J5").Value = 3.9
$ws.Range("K5").Value = 5.4
$ws.Range("F6").Value = 1.39
$ws.Range("G6").Value = 1.4
$ws.Range("H6").Value = 8.800000000000001
$ws.Range("J6").Value = 5.7
$ws.Range("K6").Value = 5.9
$ws.Range("U6").Value = 2.16
$ws.Range("X6").Value = 34
$ws.Range("Y6").Value = 42
$ws.Range("Z6").Value = 90
$ws.Range("AA6").Value = 300
$ws.Range("AB6").Value = 13
$ws.Range("AC6").Value = 14.5
$ws.Range("AD6").Value = 36
$ws.Range("AH6").Value = 24
$ws.Range("AI6").Value = 90
$ws.Range("AJ6").Value = 13
$ws.Range("AL6").Value = 29
$ws.Range("AM6").Value = 95
$ws.Range("AO6").Value = 120
$ws.Range("H7").Value = 2.12
$ws.Range("I7").Value = 2.16
$ws.Range("K7").Value = 4
$ws.Range("N7").Value = 4.8
$ws.Range("O7").Value = 1.23
$ws.Range("P7").Value = 2.32
$ws.Range("T7").Value = 1.62
$ws.Range("Y7").Value = 12.5
$ws.Range("AB7").Value = 21
$ws.Range("AD7").Value = 11
$ws.Range("AF7").Value = 34
$ws.Range("AG7").Value = 970
$ws.Range("AH7").Value = 980
$ws.Range("AN7").Value = 34
$ws.Range("AO7").Value = 55
$ws.Range("F8").Value = 1.7
$ws.Range("G8").Value = 1.72
$ws.Range("I8").Value = 6.2
$ws.Range("J8").Value = 3.95
$ws.Range("S8").Value = 3.65
$ws.Range("X8").Value = 14
$ws.Range("Y8").Value = 21
$ws.Range("AA8").Value = 210
$ws.Range("AB8").Value = 8
$ws.Range("AC8").Value = 8.800000000000001
$ws.Range("AD8").Value = 25
$ws.Range("AE8").Value = 1000
$ws.Range("AG8").Value = 11
$ws.Range("AI8").Value = 110
$ws.Range("AJ8").Value = 18.5
$ws.Range("AK8").Value = 19.5
$ws.Range("G9").Value = 2.44
$ws.Range("P9").Value = 1.74
$ws.Range("Q9").Value = 2.1
$ws.Range("F10").Value = 1.23
$ws.Range("G10").Value = 1.34
$ws.Range("H10").Value = 4.1
$ws.Range("J10").Value = 5.8
$ws.Range("K10").Value = 950
$ws.Range("P10").Value = 2.88
$ws.Range("Q10").Value = 1.41
